# Update values in columns A and B (rows 1-32) per target data, and adjust column widths.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.2562552191083114, 0.2560353230863228),
    @(-0.17669172956924406, 0.17619849683908306),
    @(-0.1264942307940391, 0.12621488001303938),
    @(-0.1182148800789804, 0.11758312223434864),
    @(-0.1145831222745608, 0.11241265536480682),
    @(-0.0224727898682886, 0.0223068921612537),
    @(-0.012306892259792424, 0.012280299450244225),
    @(-0.0022802995497350764, 0.0022610953436736914),
    @(-0.00026109539406160565, 0.0002519867901509265),
    @(0.0017480131592222392, -0.0017476697768028515),
    @(0.004747669719966652, -0.0047475344556406895),
    @(-0.04728642159202279, 0.04692971253807432),
    @(-0.01716885674490154, 0.017081341482711565),
    @(-0.009081341581530289, 0.009052818876404523),
    @(-0.008052818932472228, 0.008034400298785549),
    @(-0.006034400361719872, 0.006003453648979029),
    @(-0.004003453712960514, 0.003999999923376407),
    @(-0.0161040870501985, 0.016091395934637376),
    @(-0.012091395961800977, 0.012016627980696803),
    @(-0.008016628009970717, 0.008005653949595626),
    @(-0.0040056539791661905, 0.003999999970196733),
    @(-0.04570426608144018, 0.04549334283890438),
    @(-0.04049334288320949, 0.040097858337479586),
    @(-0.020097858478985486, 0.0199999998565481),
    @(-0.09722623444578993, 0.09710321243633935),
    @(-0.09460321248966608, 0.09444453633891747),
    @(-0.09194453639463296, 0.09100298223224446),
    @(-0.08900298229831538, 0.08835838296256959),
    @(-0.08135838306838661, 0.08117163070974165),
    @(-0.021171631149241144, 0.02102304917459863),
    @(-0.014023049288597988, 0.014000989634155303),
    @(-0.0040009897672543815, 0.003999999903951945)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Column widths: stored OOXML widths move from 15.42578125/15.42578125
# to 16.42578125/15.7109375. The engine quantizes ColumnWidth to 1/6-char
# steps, so we pick the ColumnWidth that lands on the nearest reachable width.
$ws.Columns.Item(1).ColumnWidth = 15.666666666666666
$ws.Columns.Item(2).ColumnWidth = 14.833333333333334

